# Re-order the header/template columns in the Staging.Output template.
# Columns A-E swap their labels (F and G remain untouched):
#   A2: BusinessKey        -> Output_ID
#   B2: Code                -> BusinessKey
#   C2: LongName            -> OutcomeBusinessKey
#   D2: OutcomeBusinessKey  -> Code
#   E2: Output_ID           -> LongName
#   F2: ShortName           -> ShortName        (unchanged)
#   G2: TextDescription     -> TextDescription  (unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Output_ID"
$ws.Range("B2").Value = "BusinessKey"
$ws.Range("C2").Value = "OutcomeBusinessKey"
$ws.Range("D2").Value = "Code"
$ws.Range("E2").Value = "LongName"
